$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared-string text edits -------------------------------------------------
# A4 (merged A4:A7): "Stream-lipidcane" -> "Stream-sugarcane"
$ws.Range("A4").Value = "Stream-sugarcane"

# B6: "Additional lipid extraction efficiency [%]" -> "Bagasse lipid extraction efficiency [%]"
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# --- Numeric data block (rows 4-13, columns C,D,E,F,H,I,J) -------------------
# Column D is cleared out entirely (no more data there) and a new column J is
# populated, while C/E/F/H/I get refreshed values.

$data = @{
    4  = @{ C = 0.0404635115877897;  E = -0.07889747243681094; F = -0.07028125703142579; H = -0.0767014175354384;  I = -0.08339008475211883; J = 0.1218045325129294 }
    5  = @{ C = 0.02181804545113629; E = 0.1343358583964599;    F = 0.1323993099827496;   H = 0.1364779119477987;   I = 0.08714767869196731;  J = 0.03280186734912495 }
    6  = @{ C = -0.005848646216155405; E = -0.05899947498687468; F = -0.06596114902872573; H = -0.06198154953873847; I = -0.06938873471836797; J = 0.08441035871967016 }
    7  = @{ C = 0.08461711542788571; E = 0.8631200780019502;   F = 0.8734613365334135;   H = 0.8663766594164855;   I = 0.9999099977499439;   J = -0.07982780558116795 }
    8  = @{ C = 0.9738648466211656;  E = 0.03920348008700219;  F = 0.04167554188854722;  H = 0.04057751443786095;  I = 0.0428020700517513;   J = 0.04609412241704304 }
    9  = @{ C = -0.01316282907072677; E = -0.004995124878121954; F = -0.004845121128028201; H = -0.005329633240831021; I = 0.02213005325133129; J = -0.00242407157670279 }
    10 = @{ C = 0.1096857421435536;  E = -0.02114002850071252;  F = -0.01516237905947649;  H = -0.01969699242481062;  I = -0.03204680117002925; J = -0.0129336335223706 }
    11 = @{ C = 0.05156678916972925; E = 0.1040351008775219;    F = 0.0976689417235431;   H = 0.1030765769144229;   I = 0.1248826220655516;   J = 0.002036190114259168 }
    12 = @{ C = 0.04514962874071852; E = 0.3997599939998501;    F = 0.3808070201755045;   H = 0.3954653866346658;   I = -0.09050776269406736; J = -0.01068437119276523 }
    13 = @{ C = -0.1700952523813096; E = 0.012028800720018;     F = 0.02430060751518788;  H = 0.01491787294682367;  I = 0.06003150078751969;  J = 0.08609167855977877 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]

    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").ClearContents()
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("I$row").Value = $vals.I
    $ws.Range("J$row").Value = $vals.J
}
